# Generate Report for Handoff
# Updates the localization-status report: file "e2e\b.md" now has its own
# handback xliff (no longer a content-duplicate of "e2e\a.md"), moving its
# status from "Handed back: in sync with en-US" to "Ready for handoff",
# and recording that the handback version is stale (new Error Detail).

$wb = $excel.ActiveWorkbook

$readyForHandoff = "Ready for handoff"
$newDate         = "2016-08-16 10:37:19"

# ---------------------------------------------------------------------
# Overview sheet: row 3 is "e2e\b.md"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $readyForHandoff
$overview.Range("F3").Value = $readyForHandoff
$overview.Range("G3").Value = $newDate

# ---------------------------------------------------------------------
# zh-cn sheet: row 3 is "b.md"
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $readyForHandoff
# A bare "False"/"True" literal is auto-coerced to a COM Boolean by the
# engine; prefix with an apostrophe to force plain text, then drop the
# resulting quote-prefix formatting so the cell style stays the default.
$zhcn.Range("F3").Value = "'False"
$zhcn.Range("F3").ClearFormats()
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-16 10:37:14"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/387909c2a87116c2f126e42479168bef8726a881/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd3fe7045410ef367fe8d86b7e0160b6c9d2e760/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = (40 - 5/6)

# ---------------------------------------------------------------------
# de-de sheet: row 3 is "b.md"
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $readyForHandoff
$dede.Range("F3").Value = "'False"
$dede.Range("F3").ClearFormats()
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = $newDate
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/387909c2a87116c2f126e42479168bef8726a881/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cd3fe7045410ef367fe8d86b7e0160b6c9d2e760/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = (40 - 5/6)

Write-Output "Report regenerated for handoff"
